$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.020335
$ws.Range("H2").Value = 0.061005
$ws.Range("I2").Value = 0.009804808687698561
$ws.Range("J2").Value = 0.009804808687698559
$ws.Range("M2").Value = 0.2203263333333333
$ws.Range("N2").Value = 0.660979
$ws.Range("Q2").Value = 0.004480335988333333
$ws.Range("R2").Value = 0.040323023895
$ws.Range("S2").Value = 0.009804808687698561
$ws.Range("T2").Value = 0.009804808687698559

# Row 3
$ws.Range("I3").Value = 0.1486140913768632
$ws.Range("J3").Value = 0.1486140913768632
$ws.Range("M3").Value = 0.2203263333333333
$ws.Range("N3").Value = 0.660979
$ws.Range("Q3").Value = 0.06790964343899999
$ws.Range("R3").Value = 0.611186790951
$ws.Range("S3").Value = 0.1486140913768632
$ws.Range("T3").Value = 0.1486140913768632

# Row 4
$ws.Range("H4").Value = 5.236273000000001
$ws.Range("I4").Value = 0.8415810999354383
$ws.Range("J4").Value = 0.8415810999354382
$ws.Range("M4").Value = 0.2203263333333333
$ws.Range("N4").Value = 0.660979
$ws.Range("Q4").Value = 0.3845629434741111
$ws.Range("R4").Value = 3.461066491267
$ws.Range("S4").Value = 0.8415810999354383
$ws.Range("T4").Value = 0.8415810999354382
